$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at row 1160, pushing the existing rows
# (old 1160-1240) down to (new 1163-1243).
$ws.Rows.Item(1160).Resize(3).Insert()

# Populate the 3 newly inserted rows (new weekly price records).
$ws.Range("A1160").Value = 6
$ws.Range("B1160").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1160").Value = "Metropolitana"
$ws.Range("D1160").Value = 45021
$ws.Range("E1160").Value = 13
$ws.Range("F1160").Value = 100112031
$ws.Range("G1160").Value = "Poroto verde"
$ws.Range("H1160").Value = "Magnum"
$ws.Range("I1160").Value = "Primera"
$ws.Range("J1160").Value = 580
$ws.Range("K1160").Value = 17000
$ws.Range("L1160").Value = 18000
$ws.Range("M1160").Value = 17448
$ws.Range("N1160").Value = "$/saco 25 kilos"
$ws.Range("O1160").Value = "Región Metropolitana"
$ws.Range("P1160").Value = 698
$ws.Range("Q1160").Value = 25
$ws.Range("R1160").Value = "Hortaliza"

$ws.Range("A1161").Value = 6
$ws.Range("B1161").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1161").Value = "Metropolitana"
$ws.Range("D1161").Value = 45021
$ws.Range("E1161").Value = 13
$ws.Range("F1161").Value = 100112031
$ws.Range("G1161").Value = "Poroto verde"
$ws.Range("H1161").Value = "Magnum"
$ws.Range("I1161").Value = "Segunda"
$ws.Range("J1161").Value = 270
$ws.Range("K1161").Value = 15000
$ws.Range("L1161").Value = 15000
$ws.Range("M1161").Value = 15000
$ws.Range("N1161").Value = "$/saco 25 kilos"
$ws.Range("O1161").Value = "Región Metropolitana"
$ws.Range("P1161").Value = 600
$ws.Range("Q1161").Value = 25
$ws.Range("R1161").Value = "Hortaliza"

$ws.Range("A1162").Value = 6
$ws.Range("B1162").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1162").Value = "Metropolitana"
$ws.Range("D1162").Value = 45021
$ws.Range("E1162").Value = 13
$ws.Range("F1162").Value = 100112031
$ws.Range("G1162").Value = "Poroto verde"
$ws.Range("H1162").Value = "Sin especificar"
$ws.Range("I1162").Value = "Primera"
$ws.Range("J1162").Value = 500
$ws.Range("K1162").Value = 23000
$ws.Range("L1162").Value = 25000
$ws.Range("M1162").Value = 24080
$ws.Range("N1162").Value = "$/malla 25 kilos"
$ws.Range("O1162").Value = "Región de Coquimbo"
$ws.Range("P1162").Value = 963
$ws.Range("Q1162").Value = 25
$ws.Range("R1162").Value = "Hortaliza"
